$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.782.22"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "1.655.26"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "'304.77"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").Value = "'0.3824"
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("D8").Value = "'0.3614"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'51.22"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'1.252"
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("D11").Value = "'0.08225"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "'22.71"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "'6.543"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "'7.447"
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("D16").Value = "'0.00001240"
$ws.Range("D17").Value = "1.638.97"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "'97.86"
$ws.Range("E18").Value = "  +4.16%  "
$ws.Range("D19").Value = "'0.06985"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").Value = "'6.778"
$ws.Range("E20").Value = "  +4.24%  "
$ws.Range("D21").Value = "'17.77"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'12.74"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("D24").Value = "23.785.66"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").Value = "'2.547"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("D26").Value = "'3.093"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'21.34"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "'151.30"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").Value = "'5.242"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "'135.07"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").Value = "1.826.75"
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").Value = "'6.906"
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("D33").Value = "'1.091"
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("D34").Value = "'2.129"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").Value = "'11.91"
$ws.Range("E35").Value = "  +5.69%  "
$ws.Range("D36").Value = "'0.02849"
$ws.Range("E36").Value = "  +3.99%  "
$ws.Range("D37").Value = "'0.2523"
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("D38").Value = "'6.151"
$ws.Range("E38").Value = "  +3.31%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.07236"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.08845"
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").Value = "'12.96"
$ws.Range("E41").Value = "  +8.15%  "
$ws.Range("D42").Value = "'0.7089"
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "'0.6561"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").Value = "'2.338"
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("D47").Value = "'0.9997"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "'3.970"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").Value = "'0.07986"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "'128.82"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("D51").Value = "'1.198"
$ws.Range("E51").Value = "  +1.54%  "
